$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.06842864744397358
$ws.Range("J2").Value = 0.06842864744397358
$ws.Range("M2").Value = 266.9240163333333
$ws.Range("N2").Value = 800.7720489999999
$ws.Range("O2").Value = 0.7873936103073201
$ws.Range("P2").Value = 0.78739361030732
$ws.Range("Q2").Value = 20.41265825040323
$ws.Range("R2").Value = 183.713924253629
$ws.Range("S2").Value = 0.05388027975935713
$ws.Range("T2").Value = 0.05388027975935712
$ws.Range("I3").Value = 0.06842864744397358
$ws.Range("J3").Value = 0.06842864744397358
$ws.Range("O3").Value = 0.1525285446808506
$ws.Range("P3").Value = 0.1525285446808506
$ws.Range("S3").Value = 0.0104373220091083
$ws.Range("T3").Value = 0.0104373220091083
$ws.Range("I4").Value = 0.06842864744397358
$ws.Range("J4").Value = 0.06842864744397358
$ws.Range("M4").Value = 7.729888333333332
$ws.Range("N4").Value = 23.189665
$ws.Range("O4").Value = 0.02280223700236483
$ws.Range("P4").Value = 0.02280223700236483
$ws.Range("Q4").Value = 0.5911329037738889
$ws.Range("R4").Value = 5.320196133965
$ws.Range("S4").Value = 0.001560326236768752
$ws.Range("T4").Value = 0.001560326236768752
$ws.Range("I5").Value = 0.06842864744397358
$ws.Range("J5").Value = 0.06842864744397358
$ws.Range("M5").Value = 3.652135
$ws.Range("N5").Value = 10.956405
$ws.Range("O5").Value = 0.01077335716164485
$ws.Range("P5").Value = 0.01077335716164485
$ws.Range("Q5").Value = 0.2792921546116667
$ws.Range("R5").Value = 2.513629391505
$ws.Range("S5").Value = 0.0007372062590022036
$ws.Range("T5").Value = 0.0007372062590022036
$ws.Range("I6").Value = 0.06842864744397358
$ws.Range("J6").Value = 0.06842864744397358
$ws.Range("M6").Value = 8.984181666666666
$ws.Range("N6").Value = 26.952545
$ws.Range("O6").Value = 0.02650225084781963
$ws.Range("P6").Value = 0.02650225084781963
$ws.Range("Q6").Value = 0.6870533140494445
$ws.Range("R6").Value = 6.183479826445001
$ws.Range("S6").Value = 0.001813513179737199
$ws.Range("T6").Value = 0.001813513179737199
$ws.Range("G7").Value = 0.621785
$ws.Range("H7").Value = 1.865355
$ws.Range("I7").Value = 0.5563733034589394
$ws.Range("J7").Value = 0.5563733034589394
$ws.Range("M7").Value = 266.9240163333333
$ws.Range("N7").Value = 800.7720489999999
$ws.Range("O7").Value = 0.7873936103073201
$ws.Range("P7").Value = 0.78739361030732
$ws.Range("Q7").Value = 165.9693494958217
$ws.Range("R7").Value = 1493.724145462395
$ws.Range("S7").Value = 0.4380847840891445
$ws.Range("T7").Value = 0.4380847840891445
$ws.Range("G8").Value = 0.621785
$ws.Range("H8").Value = 1.865355
$ws.Range("I8").Value = 0.5563733034589394
$ws.Range("J8").Value = 0.5563733034589394
$ws.Range("O8").Value = 0.1525285446808506
$ws.Range("P8").Value = 0.1525285446808506
$ws.Range("Q8").Value = 32.15045564104167
$ws.Range("R8").Value = 289.354100769375
$ws.Range("S8").Value = 0.0848628102758693
$ws.Range("T8").Value = 0.0848628102758693
$ws.Range("G9").Value = 0.621785
$ws.Range("H9").Value = 1.865355
$ws.Range("I9").Value = 0.5563733034589394
$ws.Range("J9").Value = 0.5563733034589394
$ws.Range("M9").Value = 7.729888333333332
$ws.Range("N9").Value = 23.189665
$ws.Range("O9").Value = 0.02280223700236483
$ws.Range("P9").Value = 0.02280223700236483
$ws.Range("Q9").Value = 4.806328617341666
$ws.Range("R9").Value = 43.256957556075
$ws.Range("S9").Value = 0.01268655592725938
$ws.Range("T9").Value = 0.01268655592725938
$ws.Range("G10").Value = 0.621785
$ws.Range("H10").Value = 1.865355
$ws.Range("I10").Value = 0.5563733034589394
$ws.Range("J10").Value = 0.5563733034589394
$ws.Range("M10").Value = 3.652135
$ws.Range("N10").Value = 10.956405
$ws.Range("O10").Value = 0.01077335716164485
$ws.Range("P10").Value = 0.01077335716164485
$ws.Range("Q10").Value = 2.270842760975
$ws.Range("R10").Value = 20.437584848775
$ws.Range("S10").Value = 0.005994008313367371
$ws.Range("T10").Value = 0.005994008313367371
$ws.Range("G11").Value = 0.621785
$ws.Range("H11").Value = 1.865355
$ws.Range("I11").Value = 0.5563733034589394
$ws.Range("J11").Value = 0.5563733034589394
$ws.Range("M11").Value = 8.984181666666666
$ws.Range("N11").Value = 26.952545
$ws.Range("O11").Value = 0.02650225084781963
$ws.Range("P11").Value = 0.02650225084781963
$ws.Range("Q11").Value = 5.586229397608333
$ws.Range("R11").Value = 50.276064578475
$ws.Range("S11").Value = 0.01474514485329888
$ws.Range("T11").Value = 0.01474514485329888
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.4193093333333334
$ws.Range("H12").Value = 1.257928
$ws.Range("I12").Value = 0.375198049097087
$ws.Range("J12").Value = 0.375198049097087
$ws.Range("M12").Value = 266.9240163333333
$ws.Range("N12").Value = 800.7720489999999
$ws.Range("O12").Value = 0.7873936103073201
$ws.Range("P12").Value = 0.78739361030732
$ws.Range("Q12").Value = 111.9237313393858
$ws.Range("R12").Value = 1007.313582054472
$ws.Range("S12").Value = 0.2954285464588185
$ws.Range("T12").Value = 0.2954285464588185
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.4193093333333334
$ws.Range("H13").Value = 1.257928
$ws.Range("I13").Value = 0.375198049097087
$ws.Range("J13").Value = 0.375198049097087
$ws.Range("O13").Value = 0.1525285446808506
$ws.Range("P13").Value = 0.1525285446808506
$ws.Range("Q13").Value = 21.68110540011111
$ws.Range("R13").Value = 195.129948601
$ws.Range("S13").Value = 0.05722841239587301
$ws.Range("T13").Value = 0.05722841239587302
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.4193093333333334
$ws.Range("H14").Value = 1.257928
$ws.Range("I14").Value = 0.375198049097087
$ws.Range("J14").Value = 0.375198049097087
$ws.Range("M14").Value = 7.729888333333332
$ws.Range("N14").Value = 23.189665
$ws.Range("O14").Value = 0.02280223700236483
$ws.Range("P14").Value = 0.02280223700236483
$ws.Range("Q14").Value = 3.241214323791111
$ws.Range("R14").Value = 29.17092891412
$ws.Range("S14").Value = 0.008555354838336691
$ws.Range("T14").Value = 0.008555354838336693
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.4193093333333334
$ws.Range("H15").Value = 1.257928
$ws.Range("I15").Value = 0.375198049097087
$ws.Range("J15").Value = 0.375198049097087
$ws.Range("M15").Value = 3.652135
$ws.Range("N15").Value = 10.956405
$ws.Range("O15").Value = 0.01077335716164485
$ws.Range("P15").Value = 0.01077335716164485
$ws.Range("Q15").Value = 1.531374292093333
$ws.Range("R15").Value = 13.78236862884
$ws.Range("S15").Value = 0.00404214258927528
$ws.Range("T15").Value = 0.004042142589275281
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.4193093333333334
$ws.Range("H16").Value = 1.257928
$ws.Range("I16").Value = 0.375198049097087
$ws.Range("J16").Value = 0.375198049097087
$ws.Range("M16").Value = 8.984181666666666
$ws.Range("N16").Value = 26.952545
$ws.Range("O16").Value = 0.02650225084781963
$ws.Range("P16").Value = 0.02650225084781963
$ws.Range("Q16").Value = 3.767151225195556
$ws.Range("R16").Value = 33.90436102676001
$ws.Range("S16").Value = 0.009943592814783544
$ws.Range("T16").Value = 0.009943592814783546
